$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 758.125
$ws.Range("I28").Value = 758.125
$ws.Range("K28").Value = 758.125
$ws.Range("M28").Value = -273.125
$ws.Range("H76").Value = 5199.1
$ws.Range("I76").Value = 4283.857
$ws.Range("J76").Value = 7334.6665
$ws.Range("K76").Value = 4283.857
$ws.Range("L76").Value = 7334.6665
$ws.Range("M76").Value = -3968.857
$ws.Range("N76").Value = -7964.6665
$ws.Range("H79").Value = 5199.1
$ws.Range("I79").Value = 4283.857
$ws.Range("J79").Value = 7334.6665
$ws.Range("K79").Value = 4283.857
$ws.Range("L79").Value = 7334.6665
$ws.Range("M79").Value = -3191.857
$ws.Range("N79").Value = -9518.666499999999
$ws.Range("H86").Value = 4626
$ws.Range("I86").Value = 4638.9165
$ws.Range("K86").Value = 4638.9165
$ws.Range("M86").Value = -3515.9165
$ws.Range("H89").Value = 4626
$ws.Range("I89").Value = 4638.9165
$ws.Range("K89").Value = 23194.5825
$ws.Range("M89").Value = -17578.5825

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13802.259
$ws.Range("I32").Value = 13993.988
$ws.Range("J32").Value = 11556.286
$ws.Range("K32").Value = 13993.988
$ws.Range("L32").Value = 11556.286
$ws.Range("M32").Value = -13706.988
$ws.Range("N32").Value = -12130.286
$ws.Range("H74").Value = 1241.375
$ws.Range("I74").Value = 1224.1333
$ws.Range("K74").Value = 1224.1333
$ws.Range("M74").Value = -350.1333
$ws.Range("H77").Value = 1241.375
$ws.Range("I77").Value = 1224.1333
$ws.Range("K77").Value = 6120.666499999999
$ws.Range("M77").Value = -1752.666499999999
$ws.Range("H97").Value = 2857.25
$ws.Range("I97").Value = 1560.75
$ws.Range("J97").Value = 3505.5
$ws.Range("K97").Value = 1560.75
$ws.Range("L97").Value = 3505.5
$ws.Range("M97").Value = -1064.75
$ws.Range("N97").Value = -4497.5
$ws.Range("H109").Value = 48110.6
$ws.Range("J109").Value = 48110.6
$ws.Range("L109").Value = 48110.6
$ws.Range("N109").Value = -50884.6
$ws.Range("H112").Value = 16831.166
$ws.Range("J112").Value = 16831.166
$ws.Range("L112").Value = 16831.166
$ws.Range("N112").Value = -19785.166
$ws.Range("H113").Value = 49999
$ws.Range("J113").Value = 49999
$ws.Range("L113").Value = 49999
$ws.Range("N113").Value = -58677
$ws.Range("H114").Value = 100000
$ws.Range("J114").Value = 100000
$ws.Range("L114").Value = 100000
$ws.Range("N114").Value = -108678
$ws.Range("H119").Value = 63950
$ws.Range("J119").Value = 63950
$ws.Range("L119").Value = 63950
$ws.Range("N119").Value = -73626

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3832.318
$ws.Range("I134").Value = 2707.0625
$ws.Range("K134").Value = 8121.1875
$ws.Range("M134").Value = -5586.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 113323.22
$ws.Range("I134").Value = 201781.8
$ws.Range("K134").Value = 605345.3999999999
$ws.Range("M134").Value = -602810.3999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 926
$ws.Range("I109").Value = 926
$ws.Range("K109").Value = 2778
$ws.Range("M109").Value = -1738
$ws.Range("H121").Value = 568.8333
$ws.Range("I121").Value = 230
$ws.Range("K121").Value = 690
$ws.Range("M121").Value = 620

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3923.5454
$ws.Range("I80").Value = 3277.5
$ws.Range("J80").Value = 4698.8
$ws.Range("K80").Value = 3277.5
$ws.Range("L80").Value = 4698.8
$ws.Range("M80").Value = -2279.5
$ws.Range("N80").Value = -6694.8
$ws.Range("H83").Value = 3923.5454
$ws.Range("I83").Value = 3277.5
$ws.Range("J83").Value = 4698.8
$ws.Range("K83").Value = 16387.5
$ws.Range("L83").Value = 23494
$ws.Range("M83").Value = -11395.5
$ws.Range("N83").Value = -33478
$ws.Range("H102").Value = 4363
$ws.Range("I102").Value = 3537.1667
$ws.Range("K102").Value = 3537.1667
$ws.Range("M102").Value = -1915.1667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9180.5
$ws.Range("I40").Value = 4764.25
$ws.Range("J40").Value = 11388.625
$ws.Range("K40").Value = 4764.25
$ws.Range("L40").Value = 11388.625
$ws.Range("M40").Value = -4628.25
$ws.Range("N40").Value = -11660.625
$ws.Range("H55").Value = 1113
$ws.Range("I55").Value = 556.6
$ws.Range("K55").Value = 556.6
$ws.Range("M55").Value = -383.6
$ws.Range("H61").Value = 8693.700000000001
$ws.Range("I61").Value = 1442.375
$ws.Range("J61").Value = 37699
$ws.Range("K61").Value = 1442.375
$ws.Range("L61").Value = 37699
$ws.Range("M61").Value = -1240.375
$ws.Range("N61").Value = -38103
$ws.Range("H68").Value = 7890.1665
$ws.Range("I68").Value = 5949
$ws.Range("K68").Value = 5949
$ws.Range("M68").Value = -5200
$ws.Range("H71").Value = 7890.1665
$ws.Range("I71").Value = 5949
$ws.Range("K71").Value = 29745
$ws.Range("M71").Value = -26001
$ws.Range("H86").Value = 60000
$ws.Range("J86").Value = 60000
$ws.Range("L86").Value = 60000
$ws.Range("N86").Value = -62372
$ws.Range("H88").Value = 784643.75
$ws.Range("J88").Value = 918199
$ws.Range("L88").Value = 918199
$ws.Range("N88").Value = -919055
$ws.Range("H89").Value = 60000
$ws.Range("J89").Value = 60000
$ws.Range("L89").Value = 180000
$ws.Range("N89").Value = -191856
$ws.Range("H91").Value = 784643.75
$ws.Range("J91").Value = 918199
$ws.Range("L91").Value = 918199
$ws.Range("N91").Value = -921163
$ws.Range("H100").Value = 2032.88
$ws.Range("I100").Value = 1057.6111
$ws.Range("J100").Value = 4540.7144
$ws.Range("K100").Value = 1057.6111
$ws.Range("L100").Value = 4540.7144
$ws.Range("M100").Value = -516.6111000000001
$ws.Range("N100").Value = -5622.7144
$ws.Range("H110").Value = 79949.5
$ws.Range("J110").Value = 79949.5
$ws.Range("L110").Value = 79949.5
$ws.Range("N110").Value = -88129.5
$ws.Range("H113").Value = 8693.700000000001
$ws.Range("I113").Value = 1442.375
$ws.Range("J113").Value = 37699
$ws.Range("K113").Value = 1442.375
$ws.Range("L113").Value = 37699
$ws.Range("M113").Value = 727.625
$ws.Range("N113").Value = -42039
$ws.Range("H133").Value = 93329.664
$ws.Range("J133").Value = 93329.664
$ws.Range("L133").Value = 93329.664
$ws.Range("N133").Value = -98389.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 73918.5
$ws.Range("J124").Value = 73918.5
$ws.Range("L124").Value = 73918.5
$ws.Range("N124").Value = -83738.5
$ws.Range("H140").Value = 69999
$ws.Range("J140").Value = 69999
$ws.Range("L140").Value = 69999
$ws.Range("N140").Value = -80359
